# Update cryptocurrency Price (D) and Volume(1h) (E) columns to refreshed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.967.64"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.639.79"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.54"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.44"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0884"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").Value = "1.872.11"
$ws.Range("D13").Value = "1.638.40"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.49"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "27.963.23"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.05"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.38"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.03"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.93"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.39"
$ws.Range("E32").Value = "  +2.69%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "1.404.57"
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.15"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.53"
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.83"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "1.781.28"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.93"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0506"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.58"
$ws.Range("E51").Value = "  -1.91%  "
